# Natmi following Dr Hou advice
# Rebuild the LR-pairs table for Cd200r1-Cd200: sending clusters M1/M2 now
# fan out to six target clusters (ECs, FAPs, M1, M2, Neutro, sCs) instead of
# four (ECs, M2/M1, Neutro, sCs), and every statistic column is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("M1","Cd200r1","Cd200","ECs",3,1,15.53304966666667,46.599149,0.5918754867343928,0.5918754867343928,2,1,70.0507755,140.101551,0.485408695871283,0.4148906539949664,1088.102175030016,6528.613050180099,0.2873015081339225,0.2455636077748213),
    @("M1","Cd200r1","Cd200","FAPs",3,1,15.53304966666667,46.599149,0.5918754867343928,0.5918754867343928,2,0.6666666666666666,0.7450486666666666,2.235146,0.005162728022151658,0.006619064379338907,11.57287794341711,104.155901490754,0.003055692160988302,0.003917661951247497),
    @("M1","Cd200r1","Cd200","M1",3,1,15.53304966666667,46.599149,0.5918754867343928,0.5918754867343928,2,0.6666666666666666,0.871287,2.613861,0.006037481860562736,0.007740574547543281,13.533744244921,121.803698204289,0.003573437514870637,0.004581456327931032),
    @("M1","Cd200r1","Cd200","M2",3,1,15.53304966666667,46.599149,0.5918754867343928,0.5918754867343928,3,1,1.056899666666667,3.170699,0.007323663231443602,0.009389570439025233,16.41687501501677,147.751875135151,0.004334696739789459,0.005557456573824926),
    @("M1","Cd200r1","Cd200","Neutro",3,1,15.53304966666667,46.599149,0.5918754867343928,0.5918754867343928,3,1,46.38387733333334,139.151632,0.321411679529899,0.4120776050862342,720.4830703512409,6484.347633161168,0.1902356942638776,0.2438986330827578),
    @("M1","Cd200r1","Cd200","sCs",3,1,15.53304966666667,46.599149,0.5918754867343928,0.5918754867343928,2,1,25.205092,50.410184,0.17465575148466,0.1492825315528919,391.5119458889026,2349.071675333416,0.1033744579209443,0.08835667102381024),
    @("M2","Cd200r1","Cd200","ECs",3,1,10.71072966666667,32.132189,0.4081245132656072,0.4081245132656072,2,1,70.0507755,140.101551,0.485408695871283,0.4148906539949664,750.2949193208566,4501.76951592514,0.1981071877373605,0.1693270462201451),
    @("M2","Cd200r1","Cd200","FAPs",3,1,10.71072966666667,32.132189,0.4081245132656072,0.4081245132656072,2,0.6666666666666666,0.7450486666666666,2.235146,0.005162728022151658,0.006619064379338907,7.98001485717711,71.82013371459401,0.002107035861163356,0.00270140242809141),
    @("M2","Cd200r1","Cd200","M1",3,1,10.71072966666667,32.132189,0.4081245132656072,0.4081245132656072,2,0.6666666666666666,0.871287,2.613861,0.006037481860562736,0.007740574547543281,9.332119519081001,83.98907567172901,0.002464044345692099,0.003159118219612249),
    @("M2","Cd200r1","Cd200","M2",3,1,10.71072966666667,32.132189,0.4081245132656072,0.4081245132656072,3,1,1.056899666666667,3.170699,0.007323663231443602,0.009389570439025233,11.32016661445678,101.881499530111,0.002988966491654144,0.003832113865200307),
    @("M2","Cd200r1","Cd200","Neutro",3,1,10.71072966666667,32.132189,0.4081245132656072,0.4081245132656072,3,1,46.38387733333334,139.151632,0.321411679529899,0.4120776050862342,496.805171009161,4471.246539082449,0.1311759852660213,0.1681789720034765),
    @("M2","Cd200r1","Cd200","sCs",3,1,10.71072966666667,32.132189,0.4081245132656072,0.4081245132656072,2,1,25.205092,50.410184,0.17465575148466,0.1492825315528919,269.9649266354627,1619.789559812776,0.07128129356371569,0.06092586052908166)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
